$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the numeric-looking columns (G:K) as Text so values like
# "155.00" or "0" are stored as literal strings, matching the source data.
$ws.Range("G16:K29").NumberFormat = "@"

$ws.Range("A16").Value = " Abu Dhabi"
$ws.Range("B16").Value = " October 30 2020"
$ws.Range("C16").Value = "Royals won by 7 wickets (with 15 balls remaining)"
$ws.Range("D16").Value = "Rajasthan Royals"
$ws.Range("E16").Value = "Kings XI Punjab"
$ws.Range("F16").Value = "Steven Smith (c)"
$ws.Range("G16").Value = "31"
$ws.Range("H16").Value = "20"
$ws.Range("I16").Value = "5"
$ws.Range("J16").Value = "0"
$ws.Range("K16").Value = "155.00"

$ws.Range("A17").Value = " Abu Dhabi"
$ws.Range("B17").Value = " October 25 2020"
$ws.Range("C17").Value = "Royals won by 8 wickets (with 10 balls remaining)"
$ws.Range("D17").Value = "Rajasthan Royals"
$ws.Range("E17").Value = "Mumbai Indians"
$ws.Range("F17").Value = "Steven Smith (c)"
$ws.Range("G17").Value = "11"
$ws.Range("H17").Value = "8"
$ws.Range("I17").Value = "1"
$ws.Range("J17").Value = "1"
$ws.Range("K17").Value = "137.50"

$ws.Range("A18").Value = " Dubai (DSC)"
$ws.Range("B18").Value = " October 14 2020"
$ws.Range("C18").Value = "Capitals won by 13 runs"
$ws.Range("D18").Value = "Rajasthan Royals"
$ws.Range("E18").Value = "Delhi Capitals"
$ws.Range("F18").Value = "Steven Smith (c)"
$ws.Range("G18").Value = "1"
$ws.Range("H18").Value = "4"
$ws.Range("I18").Value = "0"
$ws.Range("J18").Value = "0"
$ws.Range("K18").Value = "25.00"

$ws.Range("A19").Value = " Dubai (DSC)"
$ws.Range("B19").Value = " October 22 2020"
$ws.Range("C19").Value = "Sunrisers won by 8 wickets (with 11 balls remaining)"
$ws.Range("D19").Value = "Rajasthan Royals"
$ws.Range("E19").Value = "Sunrisers Hyderabad"
$ws.Range("F19").Value = "Steven Smith (c)"
$ws.Range("G19").Value = "19"
$ws.Range("H19").Value = "15"
$ws.Range("I19").Value = "2"
$ws.Range("J19").Value = "0"
$ws.Range("K19").Value = "126.66"

$ws.Range("A20").Value = " Sharjah"
$ws.Range("B20").Value = " September 22 2020"
$ws.Range("C20").Value = "Royals won by 16 runs"
$ws.Range("D20").Value = "Rajasthan Royals"
$ws.Range("E20").Value = "Chennai Super Kings"
$ws.Range("F20").Value = "Steven Smith (c)"
$ws.Range("G20").Value = "69"
$ws.Range("H20").Value = "47"
$ws.Range("I20").Value = "4"
$ws.Range("J20").Value = "4"
$ws.Range("K20").Value = "146.80"

$ws.Range("A21").Value = " Dubai (DSC)"
$ws.Range("B21").Value = " September 30 2020"
$ws.Range("C21").Value = "KKR won by 37 runs"
$ws.Range("D21").Value = "Rajasthan Royals"
$ws.Range("E21").Value = "Kolkata Knight Riders"
$ws.Range("F21").Value = "Steven Smith (c)"
$ws.Range("G21").Value = "3"
$ws.Range("H21").Value = "7"
$ws.Range("I21").Value = "0"
$ws.Range("J21").Value = "0"
$ws.Range("K21").Value = "42.85"

$ws.Range("A22").Value = " Dubai (DSC)"
$ws.Range("B22").Value = " November 01 2020"
$ws.Range("C22").Value = "KKR won by 60 runs"
$ws.Range("D22").Value = "Rajasthan Royals"
$ws.Range("E22").Value = "Kolkata Knight Riders"
$ws.Range("F22").Value = "Steven Smith (c)"
$ws.Range("G22").Value = "4"
$ws.Range("H22").Value = "4"
$ws.Range("I22").Value = "1"
$ws.Range("J22").Value = "0"
$ws.Range("K22").Value = "100.00"

$ws.Range("A23").Value = " Abu Dhabi"
$ws.Range("B23").Value = " October 19 2020"
$ws.Range("C23").Value = "Royals won by 7 wickets (with 15 balls remaining)"
$ws.Range("D23").Value = "Rajasthan Royals"
$ws.Range("E23").Value = "Chennai Super Kings"
$ws.Range("F23").Value = "Steven Smith (c)"
$ws.Range("G23").Value = "26"
$ws.Range("H23").Value = "34"
$ws.Range("I23").Value = "2"
$ws.Range("J23").Value = "0"
$ws.Range("K23").Value = "76.47"

$ws.Range("A24").Value = " Sharjah"
$ws.Range("B24").Value = " September 27 2020"
$ws.Range("C24").Value = "Royals won by 4 wickets (with 3 balls remaining)"
$ws.Range("D24").Value = "Rajasthan Royals"
$ws.Range("E24").Value = "Kings XI Punjab"
$ws.Range("F24").Value = "Steven Smith (c)"
$ws.Range("G24").Value = "50"
$ws.Range("H24").Value = "27"
$ws.Range("I24").Value = "7"
$ws.Range("J24").Value = "2"
$ws.Range("K24").Value = "185.18"

$ws.Range("A25").Value = " Dubai (DSC)"
$ws.Range("B25").Value = " October 17 2020"
$ws.Range("C25").Value = "RCB won by 7 wickets (with 2 balls remaining)"
$ws.Range("D25").Value = "Rajasthan Royals"
$ws.Range("E25").Value = "Royal Challengers Bangalore"
$ws.Range("F25").Value = "Steven Smith (c)"
$ws.Range("G25").Value = "57"
$ws.Range("H25").Value = "36"
$ws.Range("I25").Value = "6"
$ws.Range("J25").Value = "1"
$ws.Range("K25").Value = "158.33"

$ws.Range("A26").Value = " Dubai (DSC)"
$ws.Range("B26").Value = " October 11 2020"
$ws.Range("C26").Value = "Royals won by 5 wickets (with 1 ball remaining)"
$ws.Range("D26").Value = "Rajasthan Royals"
$ws.Range("E26").Value = "Sunrisers Hyderabad"
$ws.Range("F26").Value = "Steven Smith (c)"
$ws.Range("G26").Value = "5"
$ws.Range("H26").Value = "6"
$ws.Range("I26").Value = "0"
$ws.Range("J26").Value = "0"
$ws.Range("K26").Value = "83.33"

$ws.Range("A27").Value = " Abu Dhabi"
$ws.Range("B27").Value = " October 06 2020"
$ws.Range("C27").Value = "Mumbai won by 57 runs"
$ws.Range("D27").Value = "Rajasthan Royals"
$ws.Range("E27").Value = "Mumbai Indians"
$ws.Range("F27").Value = "Steven Smith (c)"
$ws.Range("G27").Value = "6"
$ws.Range("H27").Value = "7"
$ws.Range("I27").Value = "1"
$ws.Range("J27").Value = "0"
$ws.Range("K27").Value = "85.71"

$ws.Range("A28").Value = " Abu Dhabi"
$ws.Range("B28").Value = " October 03 2020"
$ws.Range("C28").Value = "RCB won by 8 wickets (with 5 balls remaining)"
$ws.Range("D28").Value = "Rajasthan Royals"
$ws.Range("E28").Value = "Royal Challengers Bangalore"
$ws.Range("F28").Value = "Steven Smith (c)"
$ws.Range("G28").Value = "5"
$ws.Range("H28").Value = "5"
$ws.Range("I28").Value = "1"
$ws.Range("J28").Value = "0"
$ws.Range("K28").Value = "100.00"

$ws.Range("A29").Value = " Sharjah"
$ws.Range("B29").Value = " October 09 2020"
$ws.Range("C29").Value = "Capitals won by 46 runs"
$ws.Range("D29").Value = "Rajasthan Royals"
$ws.Range("E29").Value = "Delhi Capitals"
$ws.Range("F29").Value = "Steven Smith (c)"
$ws.Range("G29").Value = "24"
$ws.Range("H29").Value = "17"
$ws.Range("I29").Value = "2"
$ws.Range("J29").Value = "1"
$ws.Range("K29").Value = "141.17"

# Clear the temporary Text number-format so the new cells fall back to
# the workbook default style (matching the unstyled cells elsewhere),
# while the stored values remain text.
$ws.Range("G16:K29").ClearFormats()

